# Insert three new data rows (for the 2022-01-17 / serial 44578 price report)
# at the top of the "Cebollín" weekly price block, pushing the existing
# rows 740-799 down to 743-802 (dimension grows from R799 to R802).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 740, shifting rows 740:799 down to 743:802.
$ws.Rows("740:742").Insert()

# Common (repeated) column values for this market/category block.
$mercadoId   = 9
$mercado     = "Vega Central Mapocho de Santiago"
$region      = "Metropolitana"
$fecha       = 44578
$codreg      = 13
$categoriaId = 100112037
$categoria   = "Cebollín"
$variedad    = "Sin especificar"
$unidad      = "$/paquete 36 unidades"
$origen      = "Región Metropolitana"
$kgUnidades  = 36
$clasif      = "Hortaliza"

# Row 740 - Calidad "Extra"
$ws.Cells.Item(740, 1).Value  = $mercadoId
$ws.Cells.Item(740, 2).Value  = $mercado
$ws.Cells.Item(740, 3).Value  = $region
$ws.Cells.Item(740, 4).Value  = $fecha
$ws.Cells.Item(740, 5).Value  = $codreg
$ws.Cells.Item(740, 6).Value  = $categoriaId
$ws.Cells.Item(740, 7).Value  = $categoria
$ws.Cells.Item(740, 8).Value  = $variedad
$ws.Cells.Item(740, 9).Value  = "Extra"
$ws.Cells.Item(740, 10).Value = 160
$ws.Cells.Item(740, 11).Value = 2600
$ws.Cells.Item(740, 12).Value = 2800
$ws.Cells.Item(740, 13).Value = 2700
$ws.Cells.Item(740, 14).Value = $unidad
$ws.Cells.Item(740, 15).Value = $origen
$ws.Cells.Item(740, 16).Value = 75
$ws.Cells.Item(740, 17).Value = $kgUnidades
$ws.Cells.Item(740, 18).Value = $clasif

# Row 741 - Calidad "Primera"
$ws.Cells.Item(741, 1).Value  = $mercadoId
$ws.Cells.Item(741, 2).Value  = $mercado
$ws.Cells.Item(741, 3).Value  = $region
$ws.Cells.Item(741, 4).Value  = $fecha
$ws.Cells.Item(741, 5).Value  = $codreg
$ws.Cells.Item(741, 6).Value  = $categoriaId
$ws.Cells.Item(741, 7).Value  = $categoria
$ws.Cells.Item(741, 8).Value  = $variedad
$ws.Cells.Item(741, 9).Value  = "Primera"
$ws.Cells.Item(741, 10).Value = 250
$ws.Cells.Item(741, 11).Value = 2200
$ws.Cells.Item(741, 12).Value = 2400
$ws.Cells.Item(741, 13).Value = 2300
$ws.Cells.Item(741, 14).Value = $unidad
$ws.Cells.Item(741, 15).Value = $origen
$ws.Cells.Item(741, 16).Value = 64
$ws.Cells.Item(741, 17).Value = $kgUnidades
$ws.Cells.Item(741, 18).Value = $clasif

# Row 742 - Calidad "Segunda"
$ws.Cells.Item(742, 1).Value  = $mercadoId
$ws.Cells.Item(742, 2).Value  = $mercado
$ws.Cells.Item(742, 3).Value  = $region
$ws.Cells.Item(742, 4).Value  = $fecha
$ws.Cells.Item(742, 5).Value  = $codreg
$ws.Cells.Item(742, 6).Value  = $categoriaId
$ws.Cells.Item(742, 7).Value  = $categoria
$ws.Cells.Item(742, 8).Value  = $variedad
$ws.Cells.Item(742, 9).Value  = "Segunda"
$ws.Cells.Item(742, 10).Value = 97
$ws.Cells.Item(742, 11).Value = 1900
$ws.Cells.Item(742, 12).Value = 2000
$ws.Cells.Item(742, 13).Value = 1951
$ws.Cells.Item(742, 14).Value = $unidad
$ws.Cells.Item(742, 15).Value = $origen
$ws.Cells.Item(742, 16).Value = 54
$ws.Cells.Item(742, 17).Value = $kgUnidades
$ws.Cells.Item(742, 18).Value = $clasif
